$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.249.82"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "1.891.85"
$ws.Range("E3").Value = "  +1.36%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.688"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.50%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.84"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.357"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "55.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0745"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0980"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "13.96"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.28%  "
$ws.Range("E14").Value = "  +10.66%  "
$ws.Range("D15").Value = "2.164.96"
$ws.Range("E15").Value = "  +1.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.25%  "
$ws.Range("D17").Value = "1.901.34"
$ws.Range("E17").Value = "  +1.88%  "
$ws.Range("D18").Value = "35.242.11"
$ws.Range("E18").Value = "  +1.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.82%  "
$ws.Range("D20").Value = "0.0₃0825"
$ws.Range("E20").Value = "  +1.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "243.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.91%  "
$ws.Range("E24").Value = "  +7.19%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.31%  "
$ws.Range("E30").Value = "  +1.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.35"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0597"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +17.51%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  -13.98%  "
$ws.Range("E37").Value = "  +2.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.95"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0710"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0222"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.14"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.96%  "
$ws.Range("D44").Value = "1.333.35"
$ws.Range("E44").Value = "  +4.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +14.55%  "
$ws.Range("E46").Value = "  +3.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0809"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.01%  "
$ws.Range("E48").Value = "  +0.55%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.70%  "
$ws.Range("D51").Value = "2.065.14"
$ws.Range("E51").Value = "  +0.81%  "
